# Natmi LR-pairs (Gnas-Tshr): refresh permutation-test results per Dr Hou's advice.
# Re-computes columns E:T for the existing Sending/Target-cluster combinations
# (now using 3 permutations instead of 1) and appends the new sCs-sourced rows
# (rows 14-17) that complete the 4x4 Sending-cluster x Target-cluster grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnas"
$ws.Cells.Item(2, 3).Value = "Tshr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 190.664594
$ws.Cells.Item(2, 8).Value = 571.993782
$ws.Cells.Item(2, 9).Value = 0.2001939625490346
$ws.Cells.Item(2, 10).Value = 0.2001939625490346
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.22859
$ws.Cells.Item(2, 14).Value = 0.68577
$ws.Cells.Item(2, 15).Value = 0.06730352972305123
$ws.Cells.Item(2, 16).Value = 0.06730352972305123
$ws.Cells.Item(2, 17).Value = 43.58401954246
$ws.Cells.Item(2, 18).Value = 392.25617588214
$ws.Cells.Item(2, 19).Value = 0.01347376030879436
$ws.Cells.Item(2, 20).Value = 0.01347376030879436

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnas"
$ws.Cells.Item(3, 3).Value = "Tshr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 190.664594
$ws.Cells.Item(3, 8).Value = 571.993782
$ws.Cells.Item(3, 9).Value = 0.2001939625490346
$ws.Cells.Item(3, 10).Value = 0.2001939625490346
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.180983666666667
$ws.Cells.Item(3, 14).Value = 6.542951
$ws.Cells.Item(3, 15).Value = 0.6421448840062525
$ws.Cells.Item(3, 16).Value = 0.6421448840062525
$ws.Cells.Item(3, 17).Value = 415.8363653256313
$ws.Cells.Item(3, 18).Value = 3742.527287930682
$ws.Cells.Item(3, 19).Value = 0.1285535288598019
$ws.Cells.Item(3, 20).Value = 0.1285535288598019

# Row 4: ECs -> M2
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gnas"
$ws.Cells.Item(4, 3).Value = "Tshr"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 190.664594
$ws.Cells.Item(4, 8).Value = 571.993782
$ws.Cells.Item(4, 9).Value = 0.2001939625490346
$ws.Cells.Item(4, 10).Value = 0.2001939625490346
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.173517
$ws.Cells.Item(4, 14).Value = 0.520551
$ws.Cells.Item(4, 15).Value = 0.05108844029465279
$ws.Cells.Item(4, 16).Value = 0.05108844029465278
$ws.Cells.Item(4, 17).Value = 33.083548357098
$ws.Cells.Item(4, 18).Value = 297.751935213882
$ws.Cells.Item(4, 19).Value = 0.01022759730303631
$ws.Cells.Item(4, 20).Value = 0.01022759730303631

# Row 5: ECs -> sCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Gnas"
$ws.Cells.Item(5, 3).Value = "Tshr"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 190.664594
$ws.Cells.Item(5, 8).Value = 571.993782
$ws.Cells.Item(5, 9).Value = 0.2001939625490346
$ws.Cells.Item(5, 10).Value = 0.2001939625490346
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8133136666666667
$ws.Cells.Item(5, 14).Value = 2.439941
$ws.Cells.Item(5, 15).Value = 0.2394631459760435
$ws.Cells.Item(5, 16).Value = 0.2394631459760435
$ws.Cells.Item(5, 17).Value = 155.0701200496513
$ws.Cells.Item(5, 18).Value = 1395.631080446862
$ws.Cells.Item(5, 19).Value = 0.04793907607740207
$ws.Cells.Item(5, 20).Value = 0.04793907607740206

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gnas"
$ws.Cells.Item(6, 3).Value = "Tshr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 531.1103823333333
$ws.Cells.Item(6, 8).Value = 1593.331147
$ws.Cells.Item(6, 9).Value = 0.5576551459273178
$ws.Cells.Item(6, 10).Value = 0.5576551459273177
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.22859
$ws.Cells.Item(6, 14).Value = 0.68577
$ws.Cells.Item(6, 15).Value = 0.06730352972305123
$ws.Cells.Item(6, 16).Value = 0.06730352972305123
$ws.Cells.Item(6, 17).Value = 121.4065222975767
$ws.Cells.Item(6, 18).Value = 1092.65870067819
$ws.Cells.Item(6, 19).Value = 0.03753215968913171
$ws.Cells.Item(6, 20).Value = 0.0375321596891317

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gnas"
$ws.Cells.Item(7, 3).Value = "Tshr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 531.1103823333333
$ws.Cells.Item(7, 8).Value = 1593.331147
$ws.Cells.Item(7, 9).Value = 0.5576551459273178
$ws.Cells.Item(7, 10).Value = 0.5576551459273177
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.180983666666667
$ws.Cells.Item(7, 14).Value = 6.542951
$ws.Cells.Item(7, 15).Value = 0.6421448840062525
$ws.Cells.Item(7, 16).Value = 0.6421448840062525
$ws.Cells.Item(7, 17).Value = 1158.343069066088
$ws.Cells.Item(7, 18).Value = 10425.0876215948
$ws.Cells.Item(7, 19).Value = 0.3580953989969873
$ws.Cells.Item(7, 20).Value = 0.3580953989969872

# Row 8: FAPs -> M2
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Gnas"
$ws.Cells.Item(8, 3).Value = "Tshr"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 531.1103823333333
$ws.Cells.Item(8, 8).Value = 1593.331147
$ws.Cells.Item(8, 9).Value = 0.5576551459273178
$ws.Cells.Item(8, 10).Value = 0.5576551459273177
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.173517
$ws.Cells.Item(8, 14).Value = 0.520551
$ws.Cells.Item(8, 15).Value = 0.05108844029465279
$ws.Cells.Item(8, 16).Value = 0.05108844029465278
$ws.Cells.Item(8, 17).Value = 92.156680211333
$ws.Cells.Item(8, 18).Value = 829.410121901997
$ws.Cells.Item(8, 19).Value = 0.02848973162771366
$ws.Cells.Item(8, 20).Value = 0.02848973162771365

# Row 9: FAPs -> sCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Gnas"
$ws.Cells.Item(9, 3).Value = "Tshr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 531.1103823333333
$ws.Cells.Item(9, 8).Value = 1593.331147
$ws.Cells.Item(9, 9).Value = 0.5576551459273178
$ws.Cells.Item(9, 10).Value = 0.5576551459273177
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.8133136666666667
$ws.Cells.Item(9, 14).Value = 2.439941
$ws.Cells.Item(9, 15).Value = 0.2394631459760435
$ws.Cells.Item(9, 16).Value = 0.2394631459760435
$ws.Cells.Item(9, 17).Value = 431.9593324602586
$ws.Cells.Item(9, 18).Value = 3887.633992142327
$ws.Cells.Item(9, 19).Value = 0.1335378556134851
$ws.Cells.Item(9, 20).Value = 0.1335378556134851

# Row 10: M2 -> ECs
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Gnas"
$ws.Cells.Item(10, 3).Value = "Tshr"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 79.06597733333334
$ws.Cells.Item(10, 8).Value = 237.197932
$ws.Cells.Item(10, 9).Value = 0.08301767503395074
$ws.Cells.Item(10, 10).Value = 0.08301767503395074
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.22859
$ws.Cells.Item(10, 14).Value = 0.68577
$ws.Cells.Item(10, 15).Value = 0.06730352972305123
$ws.Cells.Item(10, 16).Value = 0.06730352972305123
$ws.Cells.Item(10, 17).Value = 18.07369175862667
$ws.Cells.Item(10, 18).Value = 162.66322582764
$ws.Cells.Item(10, 19).Value = 0.005587382559186112
$ws.Cells.Item(10, 20).Value = 0.005587382559186112

# Row 11: M2 -> FAPs
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Gnas"
$ws.Cells.Item(11, 3).Value = "Tshr"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 79.06597733333334
$ws.Cells.Item(11, 8).Value = 237.197932
$ws.Cells.Item(11, 9).Value = 0.08301767503395074
$ws.Cells.Item(11, 10).Value = 0.08301767503395074
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.180983666666667
$ws.Cells.Item(11, 14).Value = 6.542951
$ws.Cells.Item(11, 15).Value = 0.6421448840062525
$ws.Cells.Item(11, 16).Value = 0.6421448840062525
$ws.Cells.Item(11, 17).Value = 172.4416051530369
$ws.Cells.Item(11, 18).Value = 1551.974446377332
$ws.Cells.Item(11, 19).Value = 0.05330937530514506
$ws.Cells.Item(11, 20).Value = 0.05330937530514506

# Row 12: M2 -> M2
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Gnas"
$ws.Cells.Item(12, 3).Value = "Tshr"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 79.06597733333334
$ws.Cells.Item(12, 8).Value = 237.197932
$ws.Cells.Item(12, 9).Value = 0.08301767503395074
$ws.Cells.Item(12, 10).Value = 0.08301767503395074
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.173517
$ws.Cells.Item(12, 14).Value = 0.520551
$ws.Cells.Item(12, 15).Value = 0.05108844029465279
$ws.Cells.Item(12, 16).Value = 0.05108844029465278
$ws.Cells.Item(12, 17).Value = 13.719291188948
$ws.Cells.Item(12, 18).Value = 123.473620700532
$ws.Cells.Item(12, 19).Value = 0.00424124353437288
$ws.Cells.Item(12, 20).Value = 0.00424124353437288

# Row 13: M2 -> sCs
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Gnas"
$ws.Cells.Item(13, 3).Value = "Tshr"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 79.06597733333334
$ws.Cells.Item(13, 8).Value = 237.197932
$ws.Cells.Item(13, 9).Value = 0.08301767503395074
$ws.Cells.Item(13, 10).Value = 0.08301767503395074
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.8133136666666667
$ws.Cells.Item(13, 14).Value = 2.439941
$ws.Cells.Item(13, 15).Value = 0.2394631459760435
$ws.Cells.Item(13, 16).Value = 0.2394631459760435
$ws.Cells.Item(13, 17).Value = 64.3054399335569
$ws.Cells.Item(13, 18).Value = 578.748959402012
$ws.Cells.Item(13, 19).Value = 0.01987967363524669
$ws.Cells.Item(13, 20).Value = 0.01987967363524669

# Row 14: sCs -> ECs
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Gnas"
$ws.Cells.Item(14, 3).Value = "Tshr"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 151.5583673333333
$ws.Cells.Item(14, 8).Value = 454.675102
$ws.Cells.Item(14, 9).Value = 0.1591332164896969
$ws.Cells.Item(14, 10).Value = 0.1591332164896969
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.22859
$ws.Cells.Item(14, 14).Value = 0.68577
$ws.Cells.Item(14, 15).Value = 0.06730352972305123
$ws.Cells.Item(14, 16).Value = 0.06730352972305123
$ws.Cells.Item(14, 17).Value = 34.64472718872666
$ws.Cells.Item(14, 18).Value = 311.80254469854
$ws.Cells.Item(14, 19).Value = 0.01071022716593906
$ws.Cells.Item(14, 20).Value = 0.01071022716593906

# Row 15: sCs -> FAPs
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Gnas"
$ws.Cells.Item(15, 3).Value = "Tshr"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 151.5583673333333
$ws.Cells.Item(15, 8).Value = 454.675102
$ws.Cells.Item(15, 9).Value = 0.1591332164896969
$ws.Cells.Item(15, 10).Value = 0.1591332164896969
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.180983666666667
$ws.Cells.Item(15, 14).Value = 6.542951
$ws.Cells.Item(15, 15).Value = 0.6421448840062525
$ws.Cells.Item(15, 16).Value = 0.6421448840062525
$ws.Cells.Item(15, 17).Value = 330.5463237006669
$ws.Cells.Item(15, 18).Value = 2974.916913306002
$ws.Cells.Item(15, 19).Value = 0.1021865808443183
$ws.Cells.Item(15, 20).Value = 0.1021865808443183

# Row 16: sCs -> M2
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Gnas"
$ws.Cells.Item(16, 3).Value = "Tshr"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 151.5583673333333
$ws.Cells.Item(16, 8).Value = 454.675102
$ws.Cells.Item(16, 9).Value = 0.1591332164896969
$ws.Cells.Item(16, 10).Value = 0.1591332164896969
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.173517
$ws.Cells.Item(16, 14).Value = 0.520551
$ws.Cells.Item(16, 15).Value = 0.05108844029465279
$ws.Cells.Item(16, 16).Value = 0.05108844029465278
$ws.Cells.Item(16, 17).Value = 26.297953224578
$ws.Cells.Item(16, 18).Value = 236.681579021202
$ws.Cells.Item(16, 19).Value = 0.008129867829529937
$ws.Cells.Item(16, 20).Value = 0.008129867829529937

# Row 17: sCs -> sCs
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Gnas"
$ws.Cells.Item(17, 3).Value = "Tshr"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 151.5583673333333
$ws.Cells.Item(17, 8).Value = 454.675102
$ws.Cells.Item(17, 9).Value = 0.1591332164896969
$ws.Cells.Item(17, 10).Value = 0.1591332164896969
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.8133136666666667
$ws.Cells.Item(17, 14).Value = 2.439941
$ws.Cells.Item(17, 15).Value = 0.2394631459760435
$ws.Cells.Item(17, 16).Value = 0.2394631459760435
$ws.Cells.Item(17, 17).Value = 123.2644914498869
$ws.Cells.Item(17, 18).Value = 1109.380423048982
$ws.Cells.Item(17, 19).Value = 0.03810654064990963
$ws.Cells.Item(17, 20).Value = 0.03810654064990962
